$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 59835.293
$ws.Range("I17").Value = 890
$ws.Range("J17").Value = 67694.664
$ws.Range("K17").Value = 2670
$ws.Range("L17").Value = 203083.992
$ws.Range("M17").Value = -2502
$ws.Range("N17").Value = -203419.992
$ws.Range("H112").Value = 17654.166
$ws.Range("I112").Value = 519
$ws.Range("J112").Value = 20579.684
$ws.Range("K112").Value = 1557
$ws.Range("L112").Value = 61739.052
$ws.Range("M112").Value = -449
$ws.Range("N112").Value = -63955.052
$ws.Range("H138").Value = 1759.4574
$ws.Range("I138").Value = 752.4286
$ws.Range("J138").Value = 2856
$ws.Range("K138").Value = 2257.2858
$ws.Range("L138").Value = 8568
$ws.Range("M138").Value = 2882.7142
$ws.Range("N138").Value = -18848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1222.4
$ws.Range("I45").Value = 1025
$ws.Range("J45").Value = 2999
$ws.Range("K45").Value = 1025
$ws.Range("L45").Value = 2999
$ws.Range("M45").Value = -648
$ws.Range("N45").Value = -3753
$ws.Range("H102").Value = 2050
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2050
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 2050
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -5294
$ws.Range("H132").Value = 1335.9574
$ws.Range("I132").Value = 945
$ws.Range("J132").Value = 2615.4546
$ws.Range("K132").Value = 2835
$ws.Range("L132").Value = 7846.3638
$ws.Range("M132").Value = -305
$ws.Range("N132").Value = -12906.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 515.1724
$ws.Range("J94").Value = 613.4167
$ws.Range("L94").Value = 613.4167
$ws.Range("N94").Value = -1515.4167
$ws.Range("H99").Value = 1622.6538
$ws.Range("I99").Value = 1045.5555
$ws.Range("J99").Value = 1928.1765
$ws.Range("K99").Value = 1045.5555
$ws.Range("L99").Value = 1928.1765
$ws.Range("M99").Value = 452.4445000000001
$ws.Range("N99").Value = -4924.1765
$ws.Range("H105").Value = 83335710
$ws.Range("I105").Value = 2350
$ws.Range("J105").Value = 125002380
$ws.Range("K105").Value = 2350
$ws.Range("L105").Value = 125002380
$ws.Range("M105").Value = -603
$ws.Range("N105").Value = -125005874
$ws.Range("H107").Value = 100000840
$ws.Range("I107").Value = 125000750
$ws.Range("J107").Value = 1213
$ws.Range("K107").Value = 125000750
$ws.Range("L107").Value = 1213
$ws.Range("M107").Value = -124998830
$ws.Range("N107").Value = -5053
$ws.Range("H134").Value = 2472759
$ws.Range("I134").Value = 778.9677
$ws.Range("J134").Value = 7946429.5
$ws.Range("K134").Value = 2336.9031
$ws.Range("L134").Value = 23839288.5
$ws.Range("M134").Value = 198.0969
$ws.Range("N134").Value = -23844358.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1755.4884
$ws.Range("I31").Value = 1154.8889
$ws.Range("J31").Value = 2769
$ws.Range("K31").Value = 1154.8889
$ws.Range("L31").Value = 2769
$ws.Range("M31").Value = -859.8888999999999
$ws.Range("N31").Value = -3359
$ws.Range("H34").Value = 1755.4884
$ws.Range("I34").Value = 1154.8889
$ws.Range("J34").Value = 2769
$ws.Range("K34").Value = 1154.8889
$ws.Range("L34").Value = 2769
$ws.Range("M34").Value = -952.8888999999999
$ws.Range("N34").Value = -3173
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H99").Value = 40001980
$ws.Range("I99").Value = 111112696
$ws.Range("J99").Value = 2206.125
$ws.Range("K99").Value = 111112696
$ws.Range("L99").Value = 2206.125
$ws.Range("M99").Value = -111111198
$ws.Range("N99").Value = -5202.125
$ws.Range("H126").Value = 40001980
$ws.Range("I126").Value = 111112696
$ws.Range("J126").Value = 2206.125
$ws.Range("K126").Value = 333338088
$ws.Range("L126").Value = 6618.375
$ws.Range("M126").Value = -333335618
$ws.Range("N126").Value = -11558.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H93").Value = 22800
$ws.Range("J93").Value = 22800
$ws.Range("L93").Value = 22800
$ws.Range("N93").Value = -26544
$ws.Range("H102").Value = 1524.3158
$ws.Range("I102").Value = 1435.25
$ws.Range("K102").Value = 1435.25
$ws.Range("M102").Value = 186.75
$ws.Range("H126").Value = 1575.1111
$ws.Range("I126").Value = 1805.8
$ws.Range("J126").Value = 1286.75
$ws.Range("K126").Value = 5417.4
$ws.Range("L126").Value = 3860.25
$ws.Range("M126").Value = -2947.4
$ws.Range("N126").Value = -8800.25
$ws.Range("H132").Value = 5042.2964
$ws.Range("I132").Value = 1089.3334
$ws.Range("K132").Value = 3268.0002
$ws.Range("M132").Value = -738.0001999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1490.0667
$ws.Range("I7").Value = 1313.2727
$ws.Range("J7").Value = 1976.25
$ws.Range("K7").Value = 1313.2727
$ws.Range("L7").Value = 1976.25
$ws.Range("M7").Value = -1201.2727
$ws.Range("N7").Value = -2200.25
$ws.Range("H40").Value = 250000000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 250000000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -250000272
$ws.Range("H68").Value = 1385.9333
$ws.Range("I68").Value = 1254
$ws.Range("K68").Value = 1254
$ws.Range("M68").Value = -505
$ws.Range("H71").Value = 1385.9333
$ws.Range("I71").Value = 1254
$ws.Range("K71").Value = 6270
$ws.Range("M71").Value = -2526
$ws.Range("H82").Value = 492.1
$ws.Range("I82").Value = 432.2143
$ws.Range("J82").Value = 631.8333
$ws.Range("K82").Value = 432.2143
$ws.Range("L82").Value = 631.8333
$ws.Range("M82").Value = -71.21429999999998
$ws.Range("N82").Value = -1353.8333
$ws.Range("H85").Value = 492.1
$ws.Range("I85").Value = 432.2143
$ws.Range("J85").Value = 631.8333
$ws.Range("K85").Value = 432.2143
$ws.Range("L85").Value = 631.8333
$ws.Range("M85").Value = 815.7857
$ws.Range("N85").Value = -3127.8333
$ws.Range("H93").Value = 1152.3684
$ws.Range("I93").Value = 1065
$ws.Range("K93").Value = 1065
$ws.Range("M93").Value = 183
$ws.Range("H95").Value = 17794.5
$ws.Range("J95").Value = 17794.5
$ws.Range("L95").Value = 17794.5
$ws.Range("N95").Value = -23286.5
$ws.Range("H100").Value = 2755.5186
$ws.Range("I100").Value = 2085.7144
$ws.Range("K100").Value = 2085.7144
$ws.Range("M100").Value = -1544.7144
$ws.Range("H126").Value = 1490.0667
$ws.Range("I126").Value = 1313.2727
$ws.Range("J126").Value = 1976.25
$ws.Range("K126").Value = 3939.8181
$ws.Range("L126").Value = 5928.75
$ws.Range("M126").Value = -1469.8181
$ws.Range("N126").Value = -10868.75
$ws.Range("H132").Value = 13432.728
$ws.Range("I132").Value = 3834.6667
$ws.Range("J132").Value = 20077.54
$ws.Range("K132").Value = 11504.0001
$ws.Range("L132").Value = 60232.62
$ws.Range("M132").Value = -8974.000100000001
$ws.Range("N132").Value = -65292.62
$ws.Range("H136").Value = 4467.467
$ws.Range("I136").Value = 4038.2917
$ws.Range("K136").Value = 12114.8751
$ws.Range("M136").Value = -9564.875100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3259
$ws.Range("I62").Value = 2967.3635
$ws.Range("J62").Value = 3579.8
$ws.Range("K62").Value = 2967.3635
$ws.Range("L62").Value = 3579.8
$ws.Range("M62").Value = -2343.3635
$ws.Range("N62").Value = -4827.8
$ws.Range("H65").Value = 3259
$ws.Range("I65").Value = 2967.3635
$ws.Range("J65").Value = 3579.8
$ws.Range("K65").Value = 14836.8175
$ws.Range("L65").Value = 17899
$ws.Range("M65").Value = -11716.8175
$ws.Range("N65").Value = -24139
$ws.Range("H97").Value = 10526
$ws.Range("J97").Value = 10526
$ws.Range("L97").Value = 10526
$ws.Range("N97").Value = -12508
$ws.Range("H122").Value = 2053.5
$ws.Range("I122").Value = 1481.1
$ws.Range("K122").Value = 4443.299999999999
$ws.Range("M122").Value = -1993.299999999999
$ws.Range("H126").Value = 896.2222
$ws.Range("I126").Value = 664.38464
$ws.Range("J126").Value = 1499
$ws.Range("K126").Value = 1993.15392
$ws.Range("L126").Value = 4497
$ws.Range("M126").Value = 476.84608
$ws.Range("N126").Value = -9437
